{"js": "// Locate the \"B1:\" paragraph, then walk two paragraphs further down (the\n// second of the three blank paragraphs sitting between \"B1:\" and \"B2:\") -\n// that blank line is where the commit adds the new sentence.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"B1:\") {\n    targetIndex = i + 2;\n    break;\n  }\n}\n\nif (targetIndex === -1 || targetIndex >= paragraphs.items.length) {\n  throw new Error(\"Could not find the blank paragraph below 'B1:'\");\n}\n\nconst target = paragraphs.items[targetIndex];\nconst inserted = target.insertText(\"Hi I am b1 branch\", Word.InsertLocation.replace);\n\n// Match the font size (10pt / half-point 20) on both the inserted run and\n// the paragraph mark itself (so the empty paragraph's own rPr also carries\n// the size, matching the authored OOXML).\ninserted.font.size = 10;\ninserted.font.sizeBidirectional = 10;\ntarget.font.size = 10;\ntarget.font.sizeBidirectional = 10;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"B1:\" paragraph, then find the empty paragraph that sits two\n# paragraphs below it (the second of the three blank paragraphs between\n# \"B1:\" and \"B2:\") - that is the line the commit adds text to.\n$paras = $d.Paragraphs\n$b1Index = $null\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $t = $paras.Item($i).Range.Text\n    $t = $t -replace \"[\\r\\a]\", \"\"\n    if ($t -eq \"B1:\") {\n        $b1Index = $i\n        break\n    }\n}\n\nif ($b1Index -eq $null) {\n    throw \"Could not find the 'B1:' paragraph\"\n}\n\n$target = $paras.Item($b1Index + 2)\n$r = $target.Range\n$r.Text = \"Hi I am b1 branch\"\n$r.Font.Size = 10\n$r.Font.SizeBi = 10\n"}
